$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-T (20 columns): Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate, Ligand average expression value,
# Ligand total expression value, Ligand derived specificity of average expression value,
# Ligand derived specificity of total expression value, Receptor-expressing cells,
# Receptor detection rate, Receptor average expression value, Receptor total expression value,
# Receptor derived specificity of average expression value,
# Receptor derived specificity of total expression value, Edge average expression weight,
# Edge total expression weight, Edge average expression derived specificity,
# Edge total expression derived specificity

$data = @(
    ,@("ECs", "Lgals1", "Cd69", "ECs", 3, 1, 16.88469033333334, 50.654071, 0.06267206522474383, 0.06267206522474383, 1, 0.3333333333333333, 13.92410533333333, 41.772316, 0.9926481335887254, 0.9926481335887255, 235.1042067220484, 2115.937860498436, 0.06221130857349283, 0.06221130857349284)
    ,@("ECs", "Lgals1", "Cd69", "FAPs", 3, 1, 16.88469033333334, 50.654071, 0.06267206522474383, 0.06267206522474383, 1, 0.3333333333333333, 0.05023433333333333, 0.150703, 0.003581200804767964, 0.003581200804767964, 0.8481911624347779, 7.633720461913001, 0.0002244412504193229, 0.0002244412504193229)
    ,@("ECs", "Lgals1", "Cd69", "sCs", 3, 1, 16.88469033333334, 50.654071, 0.06267206522474383, 0.06267206522474383, 3, 1, 0.052892, 0.158676, 0.003770665606506583, 0.003770665606506583, 0.8930650411106668, 8.037585369996, 0.0002363154008316788, 0.0002363154008316788)
    ,@("FAPs", "Lgals1", "Cd69", "ECs", 3, 1, 74.31489566666666, 222.944687, 0.2758397042791308, 0.2758397042791308, 1, 0.3333333333333333, 13.92410533333333, 41.772316, 0.9926481335887254, 0.9926481335887255, 1034.768435098343, 9312.915915885091, 0.2738117676223452, 0.2738117676223452)
    ,@("FAPs", "Lgals1", "Cd69", "FAPs", 3, 1, 74.31489566666666, 222.944687, 0.2758397042791308, 0.2758397042791308, 1, 0.3333333333333333, 0.05023433333333333, 0.150703, 0.003581200804767964, 0.003581200804767964, 3.733159240551222, 33.598433164961, 0.0009878373709513806, 0.0009878373709513806)
    ,@("FAPs", "Lgals1", "Cd69", "sCs", 3, 1, 74.31489566666666, 222.944687, 0.2758397042791308, 0.2758397042791308, 3, 1, 0.052892, 0.158676, 0.003770665606506583, 0.003770665606506583, 3.930663461601333, 35.375971154412, 0.001040099285834265, 0.001040099285834265)
    ,@("sCs", "Lgals1", "Cd69", "ECs", 3, 1, 178.2137526666667, 534.641258, 0.6614882304961253, 0.6614882304961253, 1, 0.3333333333333333, 13.92410533333333, 41.772316, 0.9926481335887254, 0.9926481335887255, 2481.467063979281, 22333.20357581352, 0.6566250573928873, 0.6566250573928873)
    ,@("sCs", "Lgals1", "Cd69", "FAPs", 3, 1, 178.2137526666667, 534.641258, 0.6614882304961253, 0.6614882304961253, 1, 0.3333333333333333, 0.05023433333333333, 0.150703, 0.003581200804767964, 0.003581200804767964, 8.952449056041555, 80.572041504374, 0.00236892218339726, 0.00236892218339726)
    ,@("sCs", "Lgals1", "Cd69", "sCs", 3, 1, 178.2137526666667, 534.641258, 0.6614882304961253, 0.6614882304961253, 3, 1, 0.052892, 0.158676, 0.003770665606506583, 0.003770665606506583, 9.426081806045334, 84.834736254408, 0.002494250919840638, 0.002494250919840638)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $row[$c]
    }
}
